# This script reflects a "Remove All" action performed on several
# day schedules: placeholder/removed workers (e.g. "Dummy Test",
# "Weekend Test", "Test Worker") are dropped, remaining shifts are
# compacted upward, and the "Full Schedule" rollup sheet is rebuilt
# to stay in sync with the per-day sheets.

$wb = $excel.ActiveWorkbook

# ---- Sunday ----
$ws = $wb.Worksheets.Item("Sunday")
$ws.Rows.Item(6).Delete()
$ws.Range("A2").Value = "12:00 PM"
$ws.Range("B2").Value = "5:00 PM"
$ws.Range("C2").Value = "Nikko Sandgren"
$ws.Range("A3").Value = "5:00 PM"
$ws.Range("B3").Value = "7:00 PM"
$ws.Range("C3").Value = "Jullian Kemp"
$ws.Range("A4").Value = "7:00 PM"
$ws.Range("B4").Value = "9:00 PM"
$ws.Range("C4").Value = "Zion Williams"
$ws.Range("A5").Value = "9:00 PM"
$ws.Range("B5").Value = "12:00 PM"
$ws.Range("C5").Value = "Olivia Schindler"

# ---- Monday ----
$ws = $wb.Worksheets.Item("Monday")
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
$ws.Range("A2").Value = "2:00 PM"
$ws.Range("B2").Value = "7:00 PM"
$ws.Range("C2").Value = "Greg Aiv"
$ws.Range("A3").Value = "7:00 PM"
$ws.Range("B3").Value = "12:00 PM"
$ws.Range("C3").Value = "Krish Chawla"

# ---- Tuesday ----
$ws = $wb.Worksheets.Item("Tuesday")
$ws.Range("A2").Value = "2:00 PM"
$ws.Range("B2").Value = "5:00 PM"
$ws.Range("C2").Value = "Sebastian Hurd"
$ws.Range("A3").Value = "5:00 PM"
$ws.Range("B3").Value = "9:00 PM"
$ws.Range("C3").Value = "Jullian Kemp"
$ws.Range("A4").Value = "9:00 PM"
$ws.Range("B4").Value = "12:00 PM"
$ws.Range("C4").Value = "Daniel Finn"

# ---- Wednesday ----
$ws = $wb.Worksheets.Item("Wednesday")
$ws.Range("A2").Value = "2:00 PM"
$ws.Range("B2").Value = "6:00 PM"
$ws.Range("C2").Value = "Alan Haim"
$ws.Range("A3").Value = "6:00 PM"
$ws.Range("B3").Value = "8:00 PM"
$ws.Range("C3").Value = "Zion Williams"
$ws.Range("A4").Value = "8:00 PM"
$ws.Range("B4").Value = "10:00 PM"
$ws.Range("C4").Value = "Gissel O Rosa"
$ws.Range("A5").Value = "10:00 PM"
$ws.Range("B5").Value = "12:00 PM"
$ws.Range("C5").Value = "Brooke Kazmierczak"

# ---- Thursday ----
$ws = $wb.Worksheets.Item("Thursday")
$ws.Rows.Item(5).Insert()
$ws.Range("A2").Value = "2:00 PM"
$ws.Range("B2").Value = "4:00 PM"
$ws.Range("C2").Value = "Daniel Senn"
$ws.Range("A3").Value = "4:00 PM"
$ws.Range("B3").Value = "8:00 PM"
$ws.Range("C3").Value = "Regenae Walkters"
$ws.Range("A4").Value = "8:00 PM"
$ws.Range("B4").Value = "10:00 PM"
$ws.Range("C4").Value = "Daniel Senn"
$ws.Range("A5").Value = "10:00 PM"
$ws.Range("B5").Value = "12:00 PM"
$ws.Range("C5").Value = "Olivia Schindler"

# ---- Friday ----
$ws = $wb.Worksheets.Item("Friday")
$ws.Range("A2").Value = "2:00 PM"
$ws.Range("B2").Value = "7:00 PM"
$ws.Range("C2").Value = "Jash Hitesh Parekh"
$ws.Range("A3").Value = "7:00 PM"
$ws.Range("B3").Value = "11:00 PM"
$ws.Range("C3").Value = "Gissel O Rosa"
$ws.Range("A4").Value = "11:00 PM"
$ws.Range("B4").Value = "12:00 PM"
$ws.Range("C4").Value = "Krish Chawla"

# ---- Saturday ----
$ws = $wb.Worksheets.Item("Saturday")
$ws.Range("A2").Value = "12:00 PM"
$ws.Range("B2").Value = "4:00 PM"
$ws.Range("C2").Value = "Tatiana Mata Diaz"
$ws.Range("A3").Value = "4:00 PM"
$ws.Range("B3").Value = "7:00 PM"
$ws.Range("C3").Value = "Daniel Finn"
$ws.Range("A4").Value = "7:00 PM"
$ws.Range("B4").Value = "11:00 PM"
$ws.Range("C4").Value = "Brooke Kazmierczak"
$ws.Range("A5").Value = "11:00 PM"
$ws.Range("B5").Value = "12:00 PM"
$ws.Range("C5").Value = "Sebastian Hurd"

# ---- Full Schedule ----
$ws = $wb.Worksheets.Item("Full Schedule")
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(26).Delete()
$ws.Range("A2").Value = "Friday"
$ws.Range("B2").Value = "2:00 PM"
$ws.Range("C2").Value = "7:00 PM"
$ws.Range("D2").Value = "Jash Hitesh Parekh"
$ws.Range("A3").Value = "Friday"
$ws.Range("B3").Value = "7:00 PM"
$ws.Range("C3").Value = "11:00 PM"
$ws.Range("D3").Value = "Gissel O Rosa"
$ws.Range("A4").Value = "Friday"
$ws.Range("B4").Value = "11:00 PM"
$ws.Range("C4").Value = "12:00 PM"
$ws.Range("D4").Value = "Krish Chawla"
$ws.Range("A5").Value = "Saturday"
$ws.Range("B5").Value = "12:00 PM"
$ws.Range("C5").Value = "4:00 PM"
$ws.Range("D5").Value = "Tatiana Mata Diaz"
$ws.Range("A6").Value = "Saturday"
$ws.Range("B6").Value = "4:00 PM"
$ws.Range("C6").Value = "7:00 PM"
$ws.Range("D6").Value = "Daniel Finn"
$ws.Range("A7").Value = "Saturday"
$ws.Range("B7").Value = "7:00 PM"
$ws.Range("C7").Value = "11:00 PM"
$ws.Range("D7").Value = "Brooke Kazmierczak"
$ws.Range("A8").Value = "Saturday"
$ws.Range("B8").Value = "11:00 PM"
$ws.Range("C8").Value = "12:00 PM"
$ws.Range("D8").Value = "Sebastian Hurd"
$ws.Range("A9").Value = "Sunday"
$ws.Range("B9").Value = "12:00 PM"
$ws.Range("C9").Value = "5:00 PM"
$ws.Range("D9").Value = "Nikko Sandgren"
$ws.Range("A10").Value = "Sunday"
$ws.Range("B10").Value = "5:00 PM"
$ws.Range("C10").Value = "7:00 PM"
$ws.Range("D10").Value = "Jullian Kemp"
$ws.Range("A11").Value = "Sunday"
$ws.Range("B11").Value = "7:00 PM"
$ws.Range("C11").Value = "9:00 PM"
$ws.Range("D11").Value = "Zion Williams"
$ws.Range("A12").Value = "Sunday"
$ws.Range("B12").Value = "9:00 PM"
$ws.Range("C12").Value = "12:00 PM"
$ws.Range("D12").Value = "Olivia Schindler"
$ws.Range("A13").Value = "Wednesday"
$ws.Range("B13").Value = "2:00 PM"
$ws.Range("C13").Value = "6:00 PM"
$ws.Range("D13").Value = "Alan Haim"
$ws.Range("A14").Value = "Wednesday"
$ws.Range("B14").Value = "6:00 PM"
$ws.Range("C14").Value = "8:00 PM"
$ws.Range("D14").Value = "Zion Williams"
$ws.Range("A15").Value = "Wednesday"
$ws.Range("B15").Value = "8:00 PM"
$ws.Range("C15").Value = "10:00 PM"
$ws.Range("D15").Value = "Gissel O Rosa"
$ws.Range("A16").Value = "Wednesday"
$ws.Range("B16").Value = "10:00 PM"
$ws.Range("C16").Value = "12:00 PM"
$ws.Range("D16").Value = "Brooke Kazmierczak"
$ws.Range("A17").Value = "Tuesday"
$ws.Range("B17").Value = "2:00 PM"
$ws.Range("C17").Value = "5:00 PM"
$ws.Range("D17").Value = "Sebastian Hurd"
$ws.Range("A18").Value = "Tuesday"
$ws.Range("B18").Value = "5:00 PM"
$ws.Range("C18").Value = "9:00 PM"
$ws.Range("D18").Value = "Jullian Kemp"
$ws.Range("A19").Value = "Tuesday"
$ws.Range("B19").Value = "9:00 PM"
$ws.Range("C19").Value = "12:00 PM"
$ws.Range("D19").Value = "Daniel Finn"
$ws.Range("A20").Value = "Monday"
$ws.Range("B20").Value = "2:00 PM"
$ws.Range("C20").Value = "7:00 PM"
$ws.Range("D20").Value = "Krish Chawla"
$ws.Range("A21").Value = "Monday"
$ws.Range("B21").Value = "7:00 PM"
$ws.Range("C21").Value = "12:00 PM"
$ws.Range("D21").Value = "Krish Chawla"
$ws.Range("A22").Value = "Thursday"
$ws.Range("B22").Value = "2:00 PM"
$ws.Range("C22").Value = "4:00 PM"
$ws.Range("D22").Value = "Daniel Senn"
$ws.Range("A23").Value = "Thursday"
$ws.Range("B23").Value = "4:00 PM"
$ws.Range("C23").Value = "8:00 PM"
$ws.Range("D23").Value = "Regenae Walkters"
$ws.Range("A24").Value = "Thursday"
$ws.Range("B24").Value = "8:00 PM"
$ws.Range("C24").Value = "10:00 PM"
$ws.Range("D24").Value = "Daniel Senn"
$ws.Range("A25").Value = "Thursday"
$ws.Range("B25").Value = "10:00 PM"
$ws.Range("C25").Value = "12:00 PM"
$ws.Range("D25").Value = "Olivia Schindler"

